$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 26: add H26 = G26/G27 ---
$ws.Range("H26").Formula = "=G26/G27"

# --- Row 29: fill D29:G29 (with existing style s=3 preserved) and H29 (style changes to s=1) ---
$ws.Range("D29").Formula = "=D26/D27"
$ws.Range("E29").Formula = "=E26/E27"
$ws.Range("F29").Formula = "=F26/F27"
$ws.Range("G29").Formula = "=G26/G27"
$ws.Range("H29").NumberFormat = "0.00"
$ws.Range("H29").Formula = "=G27/G28"

# --- Row 30: fill D30:G30 (existing style s=2 preserved); H30 stays blank ---
$ws.Range("D30").Formula = "=D26/D28"
$ws.Range("E30").Formula = "=E26/E28"
$ws.Range("F30").Formula = "=F26/F28"
$ws.Range("G30").Formula = "=G26/G28"

# --- Row 31: E31:G31 change style from s=1 ("0.00") to s=3 ("0"), remain blank; H31 gets formula, keeps s=1 ---
$ws.Range("E31:G31").NumberFormat = "0"
$ws.Range("H31").Formula = "=G25/G28"

# --- New rows 43-45 ---
$ws.Range("E43").Formula = "=E39/E40"
$ws.Range("E44").Formula = "=E39/E41"
$ws.Range("F45").Formula = "=E40/E41"

# --- New row 47 ---
$ws.Range("D47").Formula = "=D39/G39"
$ws.Range("E47").Formula = "=E39/G39"
$ws.Range("F47").Formula = "=F39/G39"

# --- New row 48 ---
$ws.Range("D48").Formula = "=D40/G40"
$ws.Range("E48").Formula = "=E40/G40"
$ws.Range("F48").Formula = "=F40/G40"

# --- New row 49 ---
$ws.Range("D49").Formula = "=D41/G41"
$ws.Range("E49").Formula = "=E41/G41"
$ws.Range("F49").Formula = "=F41/G41"

# --- View state: scroll position + selection ---
$ws.Range("F42").Select()
